# Re-curate the dimensions for the "municipio-nombre" column (column C).
# It used to be described as a measure (medida / xsd:int) and now it is
# described as a proper refArea dimension with a Municipio URI type,
# matching the newly curated dimensions mentioned in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("C3").Value = "dim"
$ws.Range("C4").Value = "URI-Municipio"
